$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A handful of the new Price values are plain decimal-looking strings
# (e.g. "200.72", "1.00", "11.80"). If assigned as-is Excel would silently
# reinterpret them as numbers (dropping the exact text, e.g. "1.00" -> 1),
# so those specific cells are switched to the Text number format first to
# preserve the literal source string, matching the original inline strings.
$numericLikeCells = @(
    "D5",
    "D6",
    "D8",
    "D10",
    "D11",
    "D12",
    "D13",
    "D14",
    "D18",
    "D20",
    "D21",
    "D22",
    "D23",
    "D24",
    "D25",
    "D26",
    "D27",
    "D28",
    "D29",
    "D30",
    "D31",
    "D32",
    "D33",
    "D34",
    "D35",
    "D36",
    "D37",
    "D39",
    "D41",
    "D42",
    "D44",
    "D45",
    "D47",
    "D48",
    "D49",
    "D50",
    "D51",
)
foreach ($addr in $numericLikeCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = "67.459.14"
$ws.Range("E2").Value = "  -0.91%  "

# Row 3
$ws.Range("D3").Value = "3.523.42"
$ws.Range("E3").Value = "  -2.39%  "

# Row 4
$ws.Range("E4").Value = "  +0.25%  "

# Row 5
$ws.Range("D5").Value = "200.72"
$ws.Range("E5").Value = "  +4.00%  "

# Row 6
$ws.Range("D6").Value = "550.48"
$ws.Range("E6").Value = "  -6.00%  "

# Row 7
$ws.Range("D7").Value = "3.513.36"
$ws.Range("E7").Value = "  -2.53%  "

# Row 8
$ws.Range("D8").Value = "0.606"
$ws.Range("E8").Value = "  -2.27%  "

# Row 9
$ws.Range("E9").Value = "  -0.02%  "

# Row 10
$ws.Range("D10").Value = "64.03"
$ws.Range("E10").Value = "  +15.82%  "

# Row 11
$ws.Range("D11").Value = "0.656"
$ws.Range("E11").Value = "  -3.19%  "

# Row 12
$ws.Range("D12").Value = "0.143"
$ws.Range("E12").Value = "  -5.36%  "

# Row 13
$ws.Range("D13").Value = "0.0000268"
$ws.Range("E13").Value = "  -6.67%  "

# Row 14
$ws.Range("D14").Value = "9.83"
$ws.Range("E14").Value = "  -1.56%  "

# Row 15
$ws.Range("D15").Value = "4.110.00"
$ws.Range("E15").Value = "  -1.96%  "

# Row 16
$ws.Range("D16").Value = "3.546.71"
$ws.Range("E16").Value = "  -2.16%  "

# Row 17
$ws.Range("E17").Value = "  -1.31%  "

# Row 18
$ws.Range("D18").Value = "18.60"
$ws.Range("E18").Value = "  +0.70%  "

# Row 19
$ws.Range("D19").Value = "67.397.32"
$ws.Range("E19").Value = "  -0.86%  "

# Row 20
$ws.Range("D20").Value = "11.82"
$ws.Range("E20").Value = "  -5.53%  "

# Row 21
$ws.Range("D21").Value = "1.03"
$ws.Range("E21").Value = "  -4.25%  "

# Row 22
$ws.Range("D22").Value = "391.54"
$ws.Range("E22").Value = "  -3.25%  "

# Row 23
$ws.Range("D23").Value = "4.03"
$ws.Range("E23").Value = "  -5.36%  "

# Row 24
$ws.Range("D24").Value = "11.95"
$ws.Range("E24").Value = "  -10.62%  "

# Row 25
$ws.Range("D25").Value = "82.54"
$ws.Range("E25").Value = "  -3.87%  "

# Row 26
$ws.Range("D26").Value = "12.26"
$ws.Range("E26").Value = "  -2.47%  "

# Row 27
$ws.Range("D27").Value = "2.81"
$ws.Range("E27").Value = "  -4.68%  "

# Row 28
$ws.Range("D28").Value = "3.74"
$ws.Range("E28").Value = "  -4.52%  "

# Row 29
$ws.Range("D29").Value = "8.85"
$ws.Range("E29").Value = "  -3.16%  "

# Row 30
$ws.Range("D30").Value = "30.89"
$ws.Range("E30").Value = "  -1.97%  "

# Row 31
$ws.Range("D31").Value = "691.31"
$ws.Range("E31").Value = "  +2.15%  "

# Row 32
$ws.Range("D32").Value = "7.18"
$ws.Range("E32").Value = "  -11.79%  "

# Row 33
$ws.Range("D33").Value = "11.80"
$ws.Range("E33").Value = "  -3.48%  "

# Row 34
$ws.Range("D34").Value = "63.91"
$ws.Range("E34").Value = "  -0.86%  "

# Row 35
$ws.Range("D35").Value = "0.111"
$ws.Range("E35").Value = "  -5.50%  "

# Row 36
$ws.Range("D36").Value = "38.86"
$ws.Range("E36").Value = "  -8.63%  "

# Row 37
$ws.Range("D37").Value = "0.404"
$ws.Range("E37").Value = "  -4.41%  "

# Row 38
$ws.Range("E38").Value = "  +0.25%  "

# Row 39
$ws.Range("D39").Value = "0.132"
$ws.Range("E39").Value = "  -2.28%  "

# Row 40
$ws.Range("B40").Value = "Maker"
$ws.Range("C40").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D40").Value = "3.078.55"
$ws.Range("E40").Value = "  -4.01%  "

# Row 41
$ws.Range("B41").Value = "ThetaToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D41").Value = "3.01"
$ws.Range("E41").Value = "  -3.35%  "

# Row 42
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  +0.16%  "

# Row 43
$ws.Range("D43").Value = "0.0₃0684"
$ws.Range("E43").Value = "  -12.81%  "

# Row 44
$ws.Range("D44").Value = "2.82"
$ws.Range("E44").Value = "  +12.06%  "

# Row 45
$ws.Range("D45").Value = "2.54"
$ws.Range("E45").Value = "  -14.10%  "

# Row 46
$ws.Range("E46").Value = "  +7.00%  "

# Row 47
$ws.Range("D47").Value = "0.0400"
$ws.Range("E47").Value = "  -5.21%  "

# Row 48
$ws.Range("D48").Value = "0.127"
$ws.Range("E48").Value = "  -3.50%  "

# Row 49
$ws.Range("D49").Value = "138.70"

# Row 50
$ws.Range("D50").Value = "8.32"
$ws.Range("E50").Value = "  -5.34%  "

# Row 51
$ws.Range("D51").Value = "2.93"
$ws.Range("E51").Value = "  -6.19%  "
